$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the new "Justifications (if any)" column (H) ---

# Header cell: copy formatting from G1 (bold header style) then set text
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Justifications (if any)"

# Data cells H2:H12: copy formatting from G2:G12 (bordered body style) then set "-"
$ws.Range("G2:G12").Copy()
$ws.Range("H2:H12").PasteSpecial(-4122)
$ws.Range("H2:H12").Value = "-"

$excel.CutCopyMode = $false

# --- Match the author's view/zoom/selection changes ---
$excel.ActiveWindow.Zoom = 70
$ws.Range("H1:H12").Select()
